$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly for the new longer labels
$ws.Columns("A").ColumnWidth = 32

# Apply the new numeric format (0.000) to the two cells that got it in this edit
$ws.Range("D3").NumberFormat = "0.000"
$ws.Range("J2").NumberFormat = "0.000"

# New row 4: "Reduce to 16 V M's"
$ws.Range("A4").Value = "Reduce to 16 V M's"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 57
$ws.Range("D4").Value = 0.067
$ws.Range("E4").Value = 73
$ws.Range("F4").Value = 0.028
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 0.029
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 0.024
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0.022

# New row 5: "Shrink Sstack and reduce to 16VM's"
$ws.Range("A5").Value = "Shrink Sstack and reduce to 16VM's"
$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 36
$ws.Range("D5").Value = 0.047
$ws.Range("E5").Value = 56
$ws.Range("F5").Value = 0.04
$ws.Range("F5").NumberFormat = "0.000"
$ws.Range("G5").Value = 54
$ws.Range("H5").Value = 0.029
$ws.Range("I5").Value = 95
$ws.Range("J5").Value = 0.025
$ws.Range("K5").Value = 45
$ws.Range("L5").Value = 0.024

# Move the active selection to the new last row, matching where the author ended up
$ws.Range("A5").Select() | Out-Null
